# Updating last week's entries and protocols:
#  - S7 now holds "Replicate # 7" (was "Replicate # 6") and is no longer the active tab
#  - a brand new sheet "S" is appended, cloned from S7's layout/protocol block,
#    holding "Replicate # 8" and this week's raw SRB readings, and becomes the active tab

$wb = $excel.ActiveWorkbook

$wsS7 = $wb.Worksheets.Item("S7")

# S7: bump the replicate label and move the selection off the old formula range
$wsS7.Range("A2").Value = "Replicate # 7"
$wsS7.Range("A2").Select()

# Clone S7 (keeps all formatting/styles) and place the clone right after it
$wsS7.Copy($null, $wsS7)
$wsNew = $wb.Worksheets.Item($wsS7.Index + 1)
$wsNew.Name = "S"

# Wipe the old S7 data (protocol header in rows 1:5 is kept) and old formulas below it
$wsNew.Range("A6:L27").Clear()

# New replicate label for the new sheet
$wsNew.Range("A2").Value = "Replicate # 8"

# NB: this PowerShell parser has no scientific-notation literal support,
# so every value below is written in plain decimal form; it serialises to
# the exact same float64 / OOXML <v> text as the source workbook.
$raw = @(
    @(0.077, 0.097, 0.113, 0.104, 0.106, 0.094, 0.091, 0.096, 0.081, 0.086, 0.079, 0.071),
    @(0.098, 0.092, 0.245, 0.255, 0.672, 0.67, 0.548, 0.507, 0.377, 0.303, 0.317, 0.081),
    @(0.095, 0.101, 0.112, 0.149, 0.277, 0.757, 0.595, 0.565, 0.412, 0.341, 0.315, 0.098),
    @(0.092, 0.129, 0.208, 0.166, 0.092, 0.773, 0.582, 0.529, 0.394, 0.317, 0.316, 0.109),
    @(0.096, 0.208, 0.154, 0.182, 0.208, 0.23, 0.632, 0.465, 0.409, 0.32, 0.329, 0.093),
    @(0.084, 0.246, 0.169, 0.339, 0.182, 0.26, 0.238, 0.458, 0.392, 0.332, 0.316, 0.088),
    @(0.09, 0.071, 0.104, 0.115, 0.104, 0.118, 0.111, 0.099, 0.091, 0.093, 0.089, 0.086),
    @(0.067, 0.091, 0.082, 0.077, 0.09, 0.091, 0.094, 0.079, 0.078, 0.075, 0.076, 0.072)
)

$startRow = 8
for ($i = 0; $i -lt $raw.Length; $i++) {
    $rowVals = $raw[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $wsNew.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

$wsNew.Range("A3").Select()

Write-Host "done"
